# Auto-generated script applying the Golem_Profits.xlsx diff to the workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 178.08333
$ws.Range("I11").Value = 178.08333
$ws.Range("K11").Value = 178.08333
$ws.Range("M11").Value = -38.08332999999999
$ws.Range("H40").Value = 2463.3635
$ws.Range("I40").Value = 1889.6
$ws.Range("K40").Value = 1889.6
$ws.Range("M40").Value = -1714.6
$ws.Range("H53").Value = 218
$ws.Range("I53").Value = 29
$ws.Range("K53").Value = 29
$ws.Range("M53").Value = 608
$ws.Range("H92").Value = 71428900
$ws.Range("I92").Value = 83333630
$ws.Range("J92").Value = 499.5
$ws.Range("K92").Value = 83333630
$ws.Range("L92").Value = 499.5
$ws.Range("M92").Value = -83332382
$ws.Range("N92").Value = -2995.5
$ws.Range("H96").Value = 750.9375
$ws.Range("I96").Value = 738
$ws.Range("J96").Value = 755.25
$ws.Range("K96").Value = 2214
$ws.Range("L96").Value = 2265.75
$ws.Range("M96").Value = -841
$ws.Range("N96").Value = -5011.75
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""
$ws.Range("H132").Value = 751
$ws.Range("I132").Value = 1001.63635
$ws.Range("J132").Value = 199.6
$ws.Range("K132").Value = 3004.90905
$ws.Range("L132").Value = 598.8
$ws.Range("M132").Value = -474.9090500000002
$ws.Range("N132").Value = -5658.8
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = ""
$ws.Range("H138").Value = 2806.6428
$ws.Range("J138").Value = 3163
$ws.Range("L138").Value = 9489
$ws.Range("N138").Value = -19769
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 11449
$ws.Range("J29").Value = 11449
$ws.Range("L29").Value = 11449
$ws.Range("N29").Value = -12065
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = ""
$ws.Range("H63").Value = 1878
$ws.Range("I63").Value = 1846
$ws.Range("K63").Value = 1846
$ws.Range("M63").Value = -1160
$ws.Range("H66").Value = 1878
$ws.Range("I66").Value = 1846
$ws.Range("K66").Value = 9230
$ws.Range("M66").Value = -5798
$ws.Range("H122").Value = 1375
$ws.Range("J122").Value = 1375
$ws.Range("L122").Value = 4125
$ws.Range("N122").Value = -9025
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2010.5
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("K58").Value = 500
$ws.Range("M58").Value = -297
$ws.Range("H86").Value = 166677170
$ws.Range("I86").Value = 250005760
$ws.Range("K86").Value = 250005760
$ws.Range("M86").Value = -250004637
$ws.Range("H89").Value = 166677170
$ws.Range("I89").Value = 250005760
$ws.Range("K89").Value = 1250028800
$ws.Range("M89").Value = -1250023184
$ws.Range("H94").Value = 1191.5555
$ws.Range("I94").Value = 1136.5
$ws.Range("J94").Value = 1235.6
$ws.Range("K94").Value = 1136.5
$ws.Range("L94").Value = 1235.6
$ws.Range("M94").Value = -685.5
$ws.Range("N94").Value = -2137.6
$ws.Range("H132").Value = 1998
$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1584.76
$ws.Range("J4").Value = 1869
$ws.Range("L4").Value = 5607
$ws.Range("N4").Value = -5831
$ws.Range("H29").Value = 31
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H55").Value = 3430.524
$ws.Range("J55").Value = 3616.2983
$ws.Range("L55").Value = 10848.8949
$ws.Range("N55").Value = -11202.8949
$ws.Range("H80").Value = 4981.75
$ws.Range("I80").Value = 4982.5
$ws.Range("J80").Value = 4981
$ws.Range("K80").Value = 14947.5
$ws.Range("L80").Value = 14943
$ws.Range("M80").Value = -14011.5
$ws.Range("N80").Value = -16815
$ws.Range("H83").Value = 4981.75
$ws.Range("I83").Value = 4982.5
$ws.Range("J83").Value = 4981
$ws.Range("K83").Value = 44842.5
$ws.Range("L83").Value = 44829
$ws.Range("M83").Value = -40162.5
$ws.Range("N83").Value = -54189
$ws.Range("H131").Value = 2196.8125
$ws.Range("I131").Value = 959.8570999999999
$ws.Range("J131").Value = 3158.889
$ws.Range("K131").Value = 2879.5713
$ws.Range("L131").Value = 9476.667000000001
$ws.Range("M131").Value = 2160.4287
$ws.Range("N131").Value = -19556.667
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 125004130
$ws.Range("J70").Value = 500000000
$ws.Range("L70").Value = 500000000
$ws.Range("N70").Value = -500000540
$ws.Range("H73").Value = 125004130
$ws.Range("J73").Value = 500000000
$ws.Range("L73").Value = 500000000
$ws.Range("N73").Value = -500001872
$ws.Range("H93").Value = 70000
$ws.Range("J93").Value = 70000
$ws.Range("L93").Value = 70000
$ws.Range("N93").Value = -73744
$ws.Range("H107").Value = 30303814
$ws.Range("J107").Value = 66668210
$ws.Range("L107").Value = 66668210
$ws.Range("N107").Value = -66672050
$ws.Range("H113").Value = 733
$ws.Range("I113").Value = 733
$ws.Range("K113").Value = 733
$ws.Range("M113").Value = 1437
$ws.Range("H122").Value = 5683.385
$ws.Range("I122").Value = 4056.8572
$ws.Range("J122").Value = 7581
$ws.Range("K122").Value = 12170.5716
$ws.Range("L122").Value = 22743
$ws.Range("M122").Value = -9720.571599999999
$ws.Range("N122").Value = -27643
$ws.Range("H126").Value = 3685.25
$ws.Range("I126").Value = 3685.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11055.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8585.75
$ws.Range("N126").Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""
$ws.Range("H22").Value = 1388.75
$ws.Range("I22").Value = 1080.3636
$ws.Range("J22").Value = 1765.6666
$ws.Range("K22").Value = 1080.3636
$ws.Range("L22").Value = 1765.6666
$ws.Range("M22").Value = -785.3635999999999
$ws.Range("N22").Value = -2355.6666
$ws.Range("H27").Value = 1388.75
$ws.Range("I27").Value = 1080.3636
$ws.Range("J27").Value = 1765.6666
$ws.Range("K27").Value = 1080.3636
$ws.Range("L27").Value = 1765.6666
$ws.Range("M27").Value = -973.3635999999999
$ws.Range("N27").Value = -1979.6666
$ws.Range("H40").Value = 1157.4
$ws.Range("I40").Value = 1157.4
$ws.Range("K40").Value = 1157.4
$ws.Range("M40").Value = -1021.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("H38").Value = 22799.8
$ws.Range("I38").Value = 22799.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 22799.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -22326.8
$ws.Range("N38").Value = ""
$ws.Range("H54").Value = 20070
$ws.Range("I54").Value = 20070
$ws.Range("K54").Value = 20070
$ws.Range("M54").Value = -19550
$ws.Range("H132").Value = 6883
$ws.Range("I132").Value = 5660.3335
$ws.Range("J132").Value = 12385
$ws.Range("K132").Value = 16981.0005
$ws.Range("L132").Value = 37155
$ws.Range("M132").Value = -14451.0005
$ws.Range("N132").Value = -42215
$ws.Range("H136").Value = 1435.1111
$ws.Range("I136").Value = 1364.5
$ws.Range("K136").Value = 4093.5
$ws.Range("M136").Value = -1543.5
